$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (45189 -> 45190) for every data row (rows 2 through 323).
$ws.Range("C2:C323").Value = 45190
